# 2019.12.07 김동욱 PSP 추가
# Fill in the new PSP time-log rows (22-27) on the "김동욱" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("김동욱")

# --- Row 22 : 2019-11-27 (Wed) 20:00-24:00 -------------------------------
$ws.Range("A22").Value = 43796
$ws.Range("B22").Value = 0.83333333333333337
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 240
$ws.Range("F22").Value = "시간 추천 안드로이드 스튜디오 코드 작성"

# --- Row 23 : 2019-11-28 (Thu) 18:00-20:00 -------------------------------
$ws.Range("A23").Value = 43797
$ws.Range("B23").Value = 0.75
$ws.Range("C23").Value = 0.83333333333333337
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 120
$ws.Range("F23").Value = "시간 추천 안드로이드 스튜디오 코드 작성"

# --- Row 24 : 2019-12-01 (Sun) 20:00-24:00 -------------------------------
$ws.Range("A24").Value = 43800
$ws.Range("B24").Value = 0.83333333333333337
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 240
$ws.Range("F24").Value = "장소 추천 안드로이드 스튜디오 코드 작성"
$ws.Range("F24").Font.Name = "돋움"

# --- Row 25 : 2019-12-02 (Mon) 14:00-22:00 -------------------------------
$ws.Range("A25").Value = 43801
$ws.Range("B25").Value = 0.58333333333333337
$ws.Range("C25").Value = 0.91666666666666663
$ws.Range("D25").Value = 120
$ws.Range("E25").Value = 360
$ws.Range("F25").Value = "장소 추천 안드로이드 스튜디오 코드 작성"
$ws.Range("F25").Font.Name = "돋움"

# --- Row 26 : 2019-12-05 (Thu) 18:00-22:00 -------------------------------
$ws.Range("A26").Value = 43804
$ws.Range("B26").Value = 0.75
$ws.Range("C26").Value = 0.91666666666666663
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 240
$ws.Range("F26").Value = "GUI 구성 코드 파악 및 추천 코드 수정"
$ws.Range("F26").Characters(5, 2).Font.Name = "맑은 고딕"
$ws.Range("F26").Characters(7, 17).Font.Name = "돋움"

# --- Row 27 : 2019-12-06 (Fri) 09:00-12:00 -------------------------------
$ws.Range("A27").Value = 43805
$ws.Range("B27").Value = 0.375
$ws.Range("C27").Value = 0.5
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 180
$ws.Range("F27").Value = "GUI 구성에 따라 추천 코드 통합"
$ws.Range("F27").Characters(5, 7).Font.Name = "맑은 고딕"
$ws.Range("F27").Characters(12, 8).Font.Name = "돋움"

# --- View state: "김동욱" becomes the active / selected sheet -----------
$ws.Activate()
$ws.Range("G15").Select()
